# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume/ranking update
# (commit: "Updated cryptos list on Wed Mar 22 18:51:55 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number + new values for the columns that changed in that row.
# B = Coin name, C = Link, D = Price, E = Volume(1h)
$updates = @(
    @{ Row=2; D='28.525.46'; E='  +1.12%  ' },
    @{ Row=3; D='1.797.28'; E='  +0.33%  ' },
    @{ Row=4; D='1.002'; E='  -0.23%  ' },
    @{ Row=5; D='330.21'; E='  -2.40%  ' },
    @{ Row=6; D='1.002'; E='  +0.13%  ' },
    @{ Row=7; D='0.4439'; E='  -1.95%  ' },
    @{ Row=8; D='0.3781'; E='  +6.01%  ' },
    @{ Row=9; D='45.56'; E='  +0.18%  ' },
    @{ Row=10; D='0.07667'; E='  +2.76%  ' },
    @{ Row=11; D='1.151'; E='  +1.21%  ' },
    @{ Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.77'; E='  +2.07%  ' },
    @{ Row=13; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.003'; E='  +0.06%  ' },
    @{ Row=14; D='6.318'; E='  +1.98%  ' },
    @{ Row=15; D='7.484'; E='  +3.61%  ' },
    @{ Row=16; D='1.784.67'; E='  -0.46%  ' },
    @{ Row=17; D='0.00001094'; E='  +1.26%  ' },
    @{ Row=18; D='0.06718'; E='  +0.67%  ' },
    @{ Row=19; D='83.85'; E='  +3.63%  ' },
    @{ Row=20; D='1.006'; E='  +0.50%  ' },
    @{ Row=21; D='17.63'; E='  +2.71%  ' },
    @{ Row=22; D='6.304'; E='  -1.04%  ' },
    @{ Row=23; D='28.510.63'; E='  +1.07%  ' },
    @{ Row=24; D='11.70'; E='  -1.08%  ' },
    @{ Row=25; D='2.429'; E='  +1.80%  ' },
    @{ Row=26; D='20.85'; E='  +2.39%  ' },
    @{ Row=27; D='2.423'; E='  +2.06%  ' },
    @{ Row=28; D='153.23'; E='  -0.27%  ' },
    @{ Row=29; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.339'; E='  +6.11%  ' },
    @{ Row=30; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='1.988.54'; E='  -0.44%  ' },
    @{ Row=31; D='131.34'; E='  -0.67%  ' },
    @{ Row=32; D='3.976'; E='  -2.29%  ' },
    @{ Row=33; D='5.942'; E='  +1.34%  ' },
    @{ Row=34; D='0.09327'; E='  -0.58%  ' },
    @{ Row=35; D='0.2265'; E='  +5.33%  ' },
    @{ Row=36; D='12.38'; E='  +2.70%  ' },
    @{ Row=37; D='0.6775'; E='  +2.41%  ' },
    @{ Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02355'; E='  -0.30%  ' },
    @{ Row=39; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06349'; E='  +2.13%  ' },
    @{ Row=40; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.281'; E='  +2.30%  ' },
    @{ Row=41; D='1.212'; E='  +0.26%  ' },
    @{ Row=42; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='8.177'; E='  +1.66%  ' },
    @{ Row=43; B='WEMIXTOKEN'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='1.451'; E='  -2.03%  ' },
    @{ Row=44; E='  +0.50%  ' },
    @{ Row=45; D='14.07'; E='  +1.72%  ' },
    @{ Row=46; D='0.6183'; E='  +2.20%  ' },
    @{ Row=47; D='3.825'; E='  -0.96%  ' },
    @{ Row=48; D='128.40'; E='  +0.16%  ' },
    @{ Row=49; D='2.048'; E='  +1.62%  ' },
    @{ Row=50; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.242'; E='  +9.84%  ' },
    @{ Row=51; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.07000'; E='  -1.12%  ' }
)


# Column D ("Price") values are digit-and-dot strings (e.g. "28.525.46", "1.002")
# that Excel's General-format type inference would otherwise coerce into numbers
# (losing formatting like trailing zeros, or misparsing the multi-dot values).
# The source workbook stores them as plain text, so force the Text number
# format while writing, then restore the cell style so no visible formatting
# change is introduced.
foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        $dcell = $ws.Range("D$row")
        $dcell.NumberFormat = "@"
        $dcell.Value = $u.D
        $dcell.Style = "Normal"
    }
    if ($u.ContainsKey('E')) {
        $ws.Range("E$row").Value = $u.E
    }
}
